$wb = $excel.ActiveWorkbook

$coverSheet = $wb.Worksheets.Item("Cover Page")
$drcSheet   = $wb.Worksheets.Item("DRC")
$histSheet  = $wb.Worksheets.Item("Temporary History")

# --- Update "Temporary History" sheet content ---
$histSheet.Range("F5").Value = "Prasanth"
$histSheet.Range("F6").Value = "Prasanth"
$histSheet.Range("F7").Value = "Prasanth"
$histSheet.Range("F8").Value = "Prasanth"
$histSheet.Range("F9").Value = "Prasanth"

$histSheet.Range("C5").Value = "1. 0"
$histSheet.Range("C6").Value = 1.1
$histSheet.Range("C7").Value = 1.2
$histSheet.Range("C8").Value = 1.3
$histSheet.Range("C9").Value = 1.4

# --- Update selections per sheet ---
$coverSheet.Activate()
$coverSheet.Range("F9").Select()

$drcSheet.Activate()
$drcSheet.Range("B40").Select()

$histSheet.Activate()
$histSheet.Range("D11").Select()
